$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value = '[''question5.py'', ''READ COMMITTED'']'
$ws.Range('B1').Value = '[''question5.py'', ''REPEATABLE READ'']'
$ws.Range('C1').Value = '[''question5.py'', ''SERIALIZABLE'']'

$ws.Range('A2').Value = 'S: 500, E: 2500, T: 1, I: READ COMMITTED'
$ws.Range('B2').Value = 'S: 500, E: 2500, T: 1, I: REPEATABLE READ'
$ws.Range('C2').Value = 'S: 500, E: 2500, T: 1, I: SERIALIZABLE'

$ws.Range('A3').Value = 'run_exchanges took 23.703131 time to finish'
$ws.Range('B3').Value = 'run_exchanges took 29.040507 time to finish'
$ws.Range('C3').Value = 'run_exchanges took 29.437218 time to finish'

$ws.Range('A4').Value = 'sum_isolation took 12.888935 time to finish'
$ws.Range('B4').Value = 'sum_isolation took 8.584639 time to finish'
$ws.Range('C4').Value = 'sum_isolation took 8.574540 time to finish'

$ws.Range('A5').Value = 'S: 500, E: 500, T: 5, I: READ COMMITTED'
$ws.Range('B5').Value = 'S: 500, E: 500, T: 5, I: REPEATABLE READ'
$ws.Range('C5').Value = 'S: 500, E: 500, T: 5, I: SERIALIZABLE'

$ws.Range('A6').Value = 'run_exchanges took 11.108206 time to finish'
$ws.Range('B6').Value = 'run_exchanges took 13.401326 time to finish'
$ws.Range('C6').Value = 'run_exchanges took 12.038902 time to finish'

$ws.Range('A7').Value = 'sum_isolation took 13.054464 time to finish'
$ws.Range('B7').Value = 'sum_isolation took 8.691852 time to finish'
$ws.Range('C7').Value = 'sum_isolation took 8.593463 time to finish'

$ws.Range('A8').Value = 'S: 500, E: 250, T: 10, I: READ COMMITTED'
$ws.Range('B8').Value = 'S: 500, E: 250, T: 10, I: REPEATABLE READ'
$ws.Range('C8').Value = 'S: 500, E: 250, T: 10, I: SERIALIZABLE'

$ws.Range('A9').Value = 'run_exchanges took 6.279267 time to finish'
$ws.Range('B9').Value = 'run_exchanges took 10.347166 time to finish'
$ws.Range('C9').Value = 'run_exchanges took 10.402270 time to finish'

$ws.Range('A10').Value = 'sum_isolation took 12.898454 time to finish'
$ws.Range('B10').Value = 'sum_isolation took 8.460448 time to finish'
$ws.Range('C10').Value = 'sum_isolation took 8.481740 time to finish'

$ws.Range('A11').Value = 'S: 500, E: 100, T: 25, I: READ COMMITTED'
$ws.Range('B11').Value = 'S: 500, E: 100, T: 25, I: REPEATABLE READ'
$ws.Range('C11').Value = 'S: 500, E: 100, T: 25, I: SERIALIZABLE'

$ws.Range('A12').Value = 'run_exchanges took 4.085643 time to finish'
$ws.Range('B12').Value = 'run_exchanges took 9.829875 time to finish'
$ws.Range('C12').Value = 'run_exchanges took 9.825487 time to finish'

$ws.Range('A13').Value = 'sum_isolation took 13.264877 time to finish'
$ws.Range('B13').Value = 'sum_isolation took 8.449010 time to finish'
$ws.Range('C13').Value = 'sum_isolation took 8.502322 time to finish'

$ws.Range('A14').Value = 'S: 500, E: 50, T: 50, I: READ COMMITTED'
$ws.Range('B14').Value = 'S: 500, E: 50, T: 50, I: REPEATABLE READ'
$ws.Range('C14').Value = 'S: 500, E: 50, T: 50, I: SERIALIZABLE'

$ws.Range('A15').Value = 'run_exchanges took 3.122220 time to finish'
$ws.Range('B15').Value = 'run_exchanges took 10.085353 time to finish'
$ws.Range('C15').Value = 'run_exchanges took 9.887863 time to finish'

$ws.Range('A16').Value = 'sum_isolation took 13.888793 time to finish'
$ws.Range('B16').Value = 'sum_isolation took 8.873881 time to finish'
$ws.Range('C16').Value = 'sum_isolation took 8.454599 time to finish'

$ws.Range("B3").Select()
